$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2 (pushing existing data down)
$ws.Rows.Item(2).Resize(2).Insert()

# The insert picks up formatting from the row above (the bold header row);
# strip it so the new data rows stay unstyled like the rest of the data.
$ws.Range("A2:C3").ClearFormats()

# Populate the two newly inserted rows with new data
$ws.Cells.Item(2,1).Value = 0.2278527319431305
$ws.Cells.Item(2,2).Value = -0.8228355050086975
$ws.Cells.Item(2,3).Value = 0.3602577745914459

$ws.Cells.Item(3,1).Value = 0.0354301854968071
$ws.Cells.Item(3,2).Value = -0.427452951669693
$ws.Cells.Item(3,3).Value = 0.4167627990245819

# The original data ran through row 22 (now shifted to row 24). Delete the
# final three rows so we end up with 20 data rows (rows 2-21) total.
$ws.Range("A22:A24").EntireRow.Delete()
